$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Qty executed upto date (numeric column C) ---
$ws.Range("C8").Value = 23
$ws.Range("C9").Value = 94
$ws.Range("C10").Value = 54
$ws.Range("C11").Value = 46
$ws.Range("C12").Value = 90
$ws.Range("C13").Value = 41
$ws.Range("C14").Value = 20
$ws.Range("C15").Value = 33
$ws.Range("C16").Value = 50
$ws.Range("C17").Value = 24

# --- Upto date Amount / Grand total cells are stored as literal text
# (e.g. "24064.00") rather than numbers, so a plain numeric-looking
# assignment would get auto-coerced into a Number and lose the
# formatted ".00" tail. Force a Text write, then restore General
# formatting so the stored style/number-format is left untouched.
function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "G9"  "24064.00"
Set-TextValue "G10" "25488.00"
Set-TextValue "G11" "30452.00"
Set-TextValue "G13" "5576.00"
Set-TextValue "G14" "460.00"

Set-TextValue "G19" "86040.00"
Set-TextValue "H19" "86040.00"
Set-TextValue "G21" "86040.00"
Set-TextValue "H21" "86040.00"
